$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308, shifting existing rows 308:346 down to 309:347.
$ws.Rows("308:308").Insert()

# Populate the newly inserted row 308 with this week's entry.
$ws.Range("A308").Value = 8
$ws.Range("B308").Value = "Terminal La Palmera de La Serena"
$ws.Range("C308").Value = "Coquimbo"
$ws.Range("D308").Value = 45131
$ws.Range("E308").Value = 4
$ws.Range("F308").Value = 100112037
$ws.Range("G308").Value = "Cebollín"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 1000
$ws.Range("K308").Value = 1000
$ws.Range("L308").Value = 1200
$ws.Range("M308").Value = 1100
$ws.Range("N308").Value = "$/paquete 6 unidades"
$ws.Range("O308").Value = "Provincia del Elquí"
$ws.Range("P308").Value = 183
$ws.Range("Q308").Value = 6
$ws.Range("R308").Value = "Hortaliza"
